# Update column F ("dSF") values on the active worksheet to match the
# repulled / recalculated data referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 2
    4  = -4
    6  = -2
    7  = 1
    8  = -4
    10 = -2
    11 = 1
    13 = -1
    14 = -2
    15 = 1
    16 = -4
    17 = -3
    18 = -3
    22 = -1
    23 = -3
    25 = 1
    28 = 3
    29 = -4
    30 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
